$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the "Third Video: Rotating Menu Button" paragraph and the
#    (empty) paragraph right before it, then fill that empty
#    paragraph with the new narrative text about the third video /
#    rotating menu button.
# ------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("Third Video: Rotating Menu Button", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Third Video: Rotating Menu Button' paragraph"
}
$thirdVideoPara = $findRange.Paragraphs(1)
$emptyPara = $thirdVideoPara.Previous()

$emptyPara.Range.InsertAfter("I also watched the third video and created the rotating menu button. Another shortcut I learned while doing this was how to select several same strings one after another by highlighting one and press Ctrl + D. I understood how to show and close (in this case) nav items. It is interesting how the example project creates the menu button itself. In the past I have used an icon. Like this I also learned how to transform single lines for example.")
# The run created directly inside a previously-empty paragraph mark
# does not automatically carry the run-level language formatting, so
# it must be (re)applied explicitly to match the rest of the document.
$emptyPara.Range.LanguageID = "en-US"

# ------------------------------------------------------------------
# 2) Insert the new ".scss split" paragraph right after it.
# ------------------------------------------------------------------
$emptyPara.Range.InsertParagraphAfter()
$scssPara = $emptyPara.Next()
$scssPara.Range.InsertAfter("Additionally, I learned how to split up .scss-files and import different ones into the main file.")

# ------------------------------------------------------------------
# 3) Insert the "26.12.2022" date paragraph right after that.
# ------------------------------------------------------------------
$scssPara.Range.InsertParagraphAfter()
$datePara = $scssPara.Next()
$datePara.Range.InsertAfter("26.12.2022")

# ------------------------------------------------------------------
# 4) Turn the old "Third Video: Rotating Menu Button" paragraph into
#    the new trailing empty paragraph by clearing its text only
#    (leaving its paragraph mark/formatting untouched). The paragraph
#    object must be re-located fresh, since the earlier structural
#    insertions invalidate previously held references.
# ------------------------------------------------------------------
$findRangeTV = $d.Content
$foundTV = $findRangeTV.Find.Execute("Third Video: Rotating Menu Button", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTV) {
    throw "Could not re-find 'Third Video: Rotating Menu Button' paragraph"
}
$thirdVideoParaFresh = $findRangeTV.Paragraphs(1)
$clearRange = $thirdVideoParaFresh.Range
$clearRange.MoveEnd(1, -1) | Out-Null
$clearRange.Delete()

# ------------------------------------------------------------------
# 5) Change the list-item text from the old "Shortcut: ..." line to
#    the new "Fourth Video: Menu Overlay & Responsiveness" heading
#    (list formatting/numbering is preserved automatically).
# ------------------------------------------------------------------
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute("Shortcut: highlight + Ctrl + D (select several same strings)", $true, $false, $false, $false, $false, $true, 1, $false, "Fourth Video: Menu Overlay & Responsiveness", 2)
if (-not $found2) {
    throw "Could not find 'Shortcut: highlight + Ctrl + D (select several same strings)' text"
}

# ------------------------------------------------------------------
# 6) Insert a new list item right after it, carrying over the same
#    list numbering/formatting.
# ------------------------------------------------------------------
$findRange3 = $d.Content
$found3 = $findRange3.Find.Execute("Fourth Video: Menu Overlay & Responsiveness", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find 'Fourth Video: Menu Overlay & Responsiveness' paragraph"
}
$fourthVideoPara = $findRange3.Paragraphs(1)
$fourthVideoPara.Range.InsertParagraphAfter()
$overlayPara = $fourthVideoPara.Next()
$overlayPara.Range.InsertAfter("First half of video: menu overlay using scss")
